# "New Pubs Edit 9-24-19"
# Remove the "Exercise 3.3" slide (slide 52, id=302 - "Please turn to the
# Exercise Manual and complete Exercise 3.3" / "Working with DataFrames and
# Series") from the deck. Its notes page is removed automatically along
# with it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(52)
$s.Delete()
